$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestSuite")

# Update the Execution Flag column (G) for rows 2-10, except row 5, from "Y" to "N"
$ws.Range("G2").Value = "N"
$ws.Range("G3").Value = "N"
$ws.Range("G4").Value = "N"
$ws.Range("G6").Value = "N"
$ws.Range("G7").Value = "N"
$ws.Range("G8").Value = "N"
$ws.Range("G9").Value = "N"
$ws.Range("G10").Value = "N"

# Update the active selection to G10
$ws.Range("G10").Select()
